$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Formula = "0.1652337620803336"
$ws.Range("H2").Formula = "25.72883082424822"
$ws.Range("I2").Formula = "0.8347300422188886"
$ws.Range("G3").Formula = "0.1655756587848228"
$ws.Range("H3").Formula = "85.96234570858464"
$ws.Range("G4").Formula = "-0.655923810589321"
$ws.Range("H4").Formula = "-6.594752182645347"
$ws.Range("G5").Formula = "-0.6371405766420055"
$ws.Range("H5").Formula = "-4.347793249060684"
$ws.Range("G6").Formula = "0.1689637756546086"
$ws.Range("H6").Formula = "-31.34561888424378"
$ws.Range("G7").Formula = "0.4254155322855189"
$ws.Range("H7").Formula = "159.6816550626114"
$ws.Range("G8").Formula = "0.08830721041578046"
$ws.Range("H8").Formula = "-46.58049551308035"
$ws.Range("G9").Formula = "0.1992121179262586"
$ws.Range("H9").Formula = "2.11160255135567"
$ws.Range("G10").Formula = "-0.1341829475738907"
$ws.Range("H10").Formula = "-134.805815476517"
$ws.Range("G11").Formula = "-0.109455549111238"
$ws.Range("H11").Formula = "7.841126043511351"
$ws.Range("G12").Formula = "0.1775623371375107"
$ws.Range("H12").Formula = "11.6538606401789"
$ws.Range("G13").Formula = "0.2117070689364375"
$ws.Range("H13").Formula = "2.939448738404239"
$ws.Range("G14").Formula = "0.1756681301008408"
$ws.Range("H14").Formula = "-7.237008806618123"
$ws.Range("G15").Formula = "0.2441129211629131"
$ws.Range("H15").Formula = "-2.317875502959454"
$ws.Range("G16").Formula = "0.01193136951809479"
$ws.Range("H16").Formula = "-67.29361631630812"
$ws.Range("G17").Formula = "0.03641696880843781"
$ws.Range("H17").Formula = "2.668066021246247"
$ws.Range("G18").Formula = "0.1181562475056396"
$ws.Range("H18").Formula = "-31.82541790109614"
$ws.Range("G19").Formula = "0.04799839224036676"
$ws.Range("H19").Formula = "-61.83119792725517"
$ws.Range("G20").Formula = "0.065762899552701"
$ws.Range("H20").Formula = "-42.64129793614874"
$ws.Range("G21").Formula = "0.06723075276597795"
$ws.Range("H21").Formula = "-33.03631577793913"
$ws.Range("G22").Formula = "0.08742614434894451"
$ws.Range("H22").Formula = "-7.188101019819934"
$ws.Range("G23").Formula = "0.08035673685045888"
$ws.Range("H23").Formula = "-25.93165652414849"
$ws.Range("G24").Formula = "-0.1815593367286216"
$ws.Range("H24").Formula = "-45.68528584878585"
$ws.Range("G25").Formula = "-0.22599331416355"
$ws.Range("H25").Formula = "-1.590992322223615"
$ws.Range("G26").Formula = "0.1809830026419857"
$ws.Range("H26").Formula = "13.83880575035217"
$ws.Range("G27").Formula = "0.1809405298450141"
$ws.Range("H27").Formula = "-9.740192252786208"
$ws.Range("G28").Formula = "0.02184480892098121"
$ws.Range("H28").Formula = "372.0247102967358"
$ws.Range("G29").Formula = "-0.0001667807417942164"
$ws.Range("H29").Formula = "-101.0845770134636"
